$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8466608812632046
$ws.Range("D2").Value = 0.006463761177446514
$ws.Range("E2").Value = 0.03949444271920299
$ws.Range("F2").Value = 3.709721482487893
$ws.Range("G2").Value = 0.002622727708114525
$ws.Range("I2").Value = 1.972153441308784
$ws.Range("J2").Value = 0.09511000923780699
$ws.Range("K2").Value = 1.644227619465426
$ws.Range("L2").Value = 0.5773587323675571
$ws.Range("M2").Value = 0.3440403399328851
$ws.Range("N2").Value = 3.698830920109231
# Row 3
$ws.Range("B3").Value = 0.831528153963859
$ws.Range("D3").Value = 0.005874905053349977
$ws.Range("E3").Value = 0.03843302258560932
$ws.Range("F3").Value = 3.70566510450827
$ws.Range("G3").Value = 0.002627270761670938
$ws.Range("I3").Value = 1.977860229400086
$ws.Range("J3").Value = 0.09350685089080812
$ws.Range("K3").Value = 1.582137032750126
$ws.Range("L3").Value = 0.5699169746221173
$ws.Range("M3").Value = 0.3387481161107466
$ws.Range("N3").Value = 3.718459829081922
# Row 4
$ws.Range("B4").Value = 0.8226512342064325
$ws.Range("D4").Value = 0.005510988232860115
$ws.Range("E4").Value = 0.03777000261431951
$ws.Range("F4").Value = 3.704672756443614
$ws.Range("G4").Value = 0.002630209688482965
$ws.Range("I4").Value = 1.981985252370819
$ws.Range("J4").Value = 0.09250253186694124
$ws.Range("K4").Value = 1.544977574468675
$ws.Range("L4").Value = 0.5656368171855917
$ws.Range("M4").Value = 0.3356690850260939
$ws.Range("N4").Value = 3.731368785670846
# Row 5
$ws.Range("B5").Value = 0.8191384072879373
$ws.Range("D5").Value = 0.005362052504068515
$ws.Range("E5").Value = 0.03749693777881546
$ws.Range("F5").Value = 3.704645147898887
$ws.Range("G5").Value = 0.00263144503192851
$ws.Range("I5").Value = 1.983822482027676
$ws.Range("J5").Value = 0.09208820326718126
$ws.Range("K5").Value = 1.530077137314606
$ws.Range("L5").Value = 0.5639653387256658
$ws.Range("M5").Value = 0.334457280438027
$ws.Range("N5").Value = 3.736844720809792
# Row 6
$ws.Range("B6").Value = 0.8185614316892611
$ws.Range("D6").Value = 0.005337282104758145
$ws.Range("E6").Value = 0.03745142063596063
$ws.Range("F6").Value = 3.704663321124713
$ws.Range("G6").Value = 0.002631652440597862
$ws.Range("I6").Value = 1.984136992252658
$ws.Range("J6").Value = 0.09201909757176097
$ws.Range("K6").Value = 1.52761756349264
$ws.Range("L6").Value = 0.5636921838113693
$ws.Range("M6").Value = 0.3342586559887764
$ws.Range("N6").Value = 3.737767007196503
# Row 7
$ws.Range("B7").Value = 0.8226034350796567
$ws.Range("D7").Value = 0.005508982268430884
$ws.Range("E7").Value = 0.03776633167015042
$ws.Range("F7").Value = 3.704670858509942
$ws.Range("G7").Value = 0.0026302261958791
$ws.Range("I7").Value = 1.982009397137951
$ws.Range("J7").Value = 0.09249696462532242
$ws.Range("K7").Value = 1.544775640937871
$ws.Range("L7").Value = 0.5656139805721949
$ws.Range("M7").Value = 0.3356525683044005
$ws.Range("N7").Value = 3.731441763793399
# Row 8
$ws.Range("B8").Value = 0.8413572474229341
$ws.Range("D8").Value = 0.006261186021649934
$ws.Range("E8").Value = 0.0391307889839112
$ws.Range("F8").Value = 3.708011878501651
$ws.Range("G8").Value = 0.002624263202253002
$ws.Range("I8").Value = 1.973992295859958
$ws.Range("J8").Value = 0.0945613624878554
$ws.Range("K8").Value = 1.622618495376088
$ws.Range("L8").Value = 0.5747328172732864
$ws.Range("M8").Value = 0.3421802524385633
$ws.Range("N8").Value = 3.705421154036245
# Row 9
$ws.Range("B9").Value = 0.8814119754782439
$ws.Range("D9").Value = 0.007719596626749592
$ws.Range("E9").Value = 0.04171852375276508
$ws.Range("F9").Value = 3.726455299446769
$ws.Range("G9").Value = 0.002613750302334341
$ws.Range("I9").Value = 1.963195546226267
$ws.Range("J9").Value = 0.09845300290917436
$ws.Range("K9").Value = 1.782938254467524
$ws.Range("L9").Value = 0.5949096345426454
$ws.Range("M9").Value = 0.3563310959780424
$ws.Range("N9").Value = 3.661192484055903
# Row 10
$ws.Range("B10").Value = 0.9128279269003485
$ws.Range("D10").Value = 0.008783840036880264
$ws.Range("E10").Value = 0.04356862440401166
$ws.Range("F10").Value = 3.74726669226537
$ws.Range("G10").Value = 0.002606738454707907
$ws.Range("I10").Value = 1.958262980676501
$ws.Range("J10").Value = 0.1012194345157198
$ws.Range("K10").Value = 1.905441157707855
$ws.Range("L10").Value = 0.611136238698208
$ws.Range("M10").Value = 0.3675495969078
$ws.Range("N10").Value = 3.632841482287205
# Row 11
$ws.Range("B11").Value = 0.9275494708644771
$ws.Range("D11").Value = 0.009267053524954605
$ws.Range("E11").Value = 0.0443997362594768
$ws.Range("F11").Value = 3.758314051646167
$ws.Range("G11").Value = 0.00260370154507461
$ws.Range("I11").Value = 1.9566700341542
$ws.Range("J11").Value = 0.1024584145074101
$ws.Range("K11").Value = 1.962204617982025
$ws.Range("L11").Value = 0.6188237225423592
$ws.Range("M11").Value = 0.3728315242343285
$ws.Range("N11").Value = 3.62084371839552
# Row 12
$ws.Range("B12").Value = 0.9331857545505784
$ws.Range("D12").Value = 0.009449952901633196
$ws.Range("E12").Value = 0.04471298757166586
$ws.Range("F12").Value = 3.762724768523029
$ws.Range("G12").Value = 0.002602573396360965
$ws.Range("I12").Value = 1.956160381249859
$ws.Range("J12").Value = 0.1029248264338065
$ws.Range("K12").Value = 1.983848959577529
$ws.Range("L12").Value = 0.6217787978199709
$ws.Range("M12").Value = 0.3748572889537272
$ws.Range("N12").Value = 3.616429826136823
# Row 13
$ws.Range("B13").Value = 0.9319691473927492
$ws.Range("D13").Value = 0.009410565162550455
$ws.Range("E13").Value = 0.04464558818261111
$ws.Range("F13").Value = 3.761764729888966
$ws.Range("G13").Value = 0.002602815393006184
$ws.Range("I13").Value = 1.956265983523373
$ws.Range("J13").Value = 0.1028244985412776
$ws.Range("K13").Value = 1.979180821500279
$ws.Range("L13").Value = 0.6211404128890479
$ws.Range("M13").Value = 0.3744198661385667
$ws.Range("N13").Value = 3.617374680201948
# Row 14
$ws.Range("B14").Value = 0.9280119394829853
$ws.Range("D14").Value = 0.009282102151598082
$ws.Range("E14").Value = 0.04442553692803841
$ws.Range("F14").Value = 3.758672367526756
$ws.Range("G14").Value = 0.002603608294073221
$ws.Range("I14").Value = 1.956626229800271
$ws.Range("J14").Value = 0.1024968416470955
$ws.Range("K14").Value = 1.963982318390549
$ws.Range("L14").Value = 0.6190659566849206
$ws.Range("M14").Value = 0.3729976721553072
$ws.Range("N14").Value = 3.620477990913372
# Row 15
$ws.Range("B15").Value = 0.925596042200624
$ws.Range("D15").Value = 0.009203405480068483
$ws.Range("E15").Value = 0.04429055853338326
$ws.Range("F15").Value = 3.756807813311568
$ws.Range("G15").Value = 0.00260409681274688
$ws.Range("I15").Value = 1.956859074426802
$ws.Range("J15").Value = 0.1022957836816545
$ws.Range("K15").Value = 1.954692240461156
$ws.Range("L15").Value = 0.6178010214859597
$ws.Range("M15").Value = 0.3721298709847076
$ws.Range("N15").Value = 3.622395712841438
# Row 16
$ws.Range("B16").Value = 0.9118744702843173
$ws.Range("D16").Value = 0.008752246611187786
$ws.Range("E16").Value = 0.04351410111087617
$ws.Range("F16").Value = 3.746576526216529
$ws.Range("G16").Value = 0.002606939989385349
$ws.Range("I16").Value = 1.958380178023887
$ws.Range("J16").Value = 0.1011380755127824
$ws.Range("K16").Value = 1.901752404541185
$ws.Range("L16").Value = 0.6106399994827427
$ws.Range("M16").Value = 0.3672079983042806
$ws.Range("N16").Value = 3.633643667112565
# Row 17
$ws.Range("B17").Value = 0.9035667058525974
$ws.Range("D17").Value = 0.008475275066334831
$ws.Range("E17").Value = 0.04303510544073141
$ws.Range("F17").Value = 3.740704744798023
$ws.Range("G17").Value = 0.002608723246013072
$ws.Range("I17").Value = 1.959480009927674
$ws.Range("J17").Value = 0.1004228963845613
$ws.Range("K17").Value = 1.869541109687816
$ws.Range("L17").Value = 0.6063252979392928
$ws.Range("M17").Value = 0.3642342792804669
$ws.Range("N17").Value = 3.640774304184035
# Row 18
$ws.Range("B18").Value = 0.898828817381883
$ws.Range("D18").Value = 0.008315879183374619
$ws.Range("E18").Value = 0.04275860840850854
$ws.Range("F18").Value = 3.737476176553344
$ws.Range("G18").Value = 0.002609763318570657
$ws.Range("I18").Value = 1.960173872257101
$ws.Range("J18").Value = 0.1000097080422222
$ws.Range("K18").Value = 1.851111549272275
$ws.Range("L18").Value = 0.6038723901926772
$ws.Range("M18").Value = 0.362540686914322
$ws.Range("N18").Value = 3.644960293967159
# Row 19
$ws.Range("B19").Value = 0.8972316184667477
$ws.Range("D19").Value = 0.008261893686579924
$ws.Range("E19").Value = 0.04266481982751635
$ws.Range("F19").Value = 3.73640857964206
$ws.Range("G19").Value = 0.002610117944038315
$ws.Range("I19").Value = 1.960419326155268
$ws.Range("J19").Value = 0.09986949301718084
$ws.Range("K19").Value = 1.844888359655727
$ws.Range("L19").Value = 0.6030468232566903
$ws.Range("M19").Value = 0.3619701551543741
$ws.Range("N19").Value = 3.646392131820704
# Row 20
$ws.Range("B20").Value = 0.9044468904268115
$ws.Range("D20").Value = 0.008504768038154253
$ws.Range("E20").Value = 0.0430861976744854
$ws.Range("F20").Value = 3.741314412830832
$ws.Range("G20").Value = 0.002608531926871069
$ws.Range("I20").Value = 1.959356590228452
$ws.Range("J20").Value = 0.1004992181835007
$ws.Range("K20").Value = 1.872959961938932
$ws.Range("L20").Value = 0.6067816251712799
$ws.Range("M20").Value = 0.3645490972318512
$ws.Range("N20").Value = 3.640006475602391
# Row 21
$ws.Range("B21").Value = 0.9291725993188891
$ws.Range("D21").Value = 0.009319836716414898
$ws.Range("E21").Value = 0.04449021097773809
$ws.Range("F21").Value = 3.75957449904918
$ws.Range("G21").Value = 0.002603374807192865
$ws.Range("I21").Value = 1.956517877793807
$ws.Range("J21").Value = 0.1025931570938781
$ws.Range("K21").Value = 1.968442434608278
$ws.Range("L21").Value = 0.6196740808080108
$ws.Range("M21").Value = 0.3734147103758332
$ws.Range("N21").Value = 3.619562960699042
# Row 22
$ws.Range("B22").Value = 0.9456909204958492
$ws.Range("D22").Value = 0.009852069875147862
$ws.Range("E22").Value = 0.0453992487370396
$ws.Range("F22").Value = 3.772833576579231
$ws.Range("G22").Value = 0.002600131715355676
$ws.Range("I22").Value = 1.95520794680489
$ws.Range("J22").Value = 0.1039455810327929
$ws.Range("K22").Value = 2.031715816461201
$ws.Range("L22").Value = 0.6283564669515158
$ws.Range("M22").Value = 0.3793581664683359
$ws.Range("N22").Value = 3.606956232360503
# Row 23
$ws.Range("B23").Value = 0.9368420702443814
$ws.Range("D23").Value = 0.009568033023519718
$ws.Range("E23").Value = 0.04491484923384093
$ws.Range("F23").Value = 3.765635680333403
$ws.Range("G23").Value = 0.002601850995410644
$ws.Range("I23").Value = 1.955857195497778
$ws.Range("J23").Value = 0.1032252259345015
$ws.Range("K23").Value = 1.997865930089745
$ws.Range("L23").Value = 0.6236990505746434
$ws.Range("M23").Value = 0.376172395609828
$ws.Range("N23").Value = 3.613615632056408
# Row 24
$ws.Range("B24").Value = 0.904048839483437
$ws.Range("D24").Value = 0.00849143477130454
$ws.Range("E24").Value = 0.04306310235539712
$ws.Range("F24").Value = 3.741038323412496
$ws.Range("G24").Value = 0.002608618376070895
$ws.Range("I24").Value = 1.959412196533421
$ws.Range("J24").Value = 0.1004647194011135
$ws.Range("K24").Value = 1.871414021190503
$ws.Range("L24").Value = 0.6065752334364163
$ws.Range("M24").Value = 0.3644067180740294
$ws.Range("N24").Value = 3.640353341488989
# Row 25
$ws.Range("B25").Value = 0.870226260544797
$ws.Range("D25").Value = 0.007326542874832143
$ws.Range("E25").Value = 0.04102763065728787
$ws.Range("F25").Value = 3.720191617287682
$ws.Range("G25").Value = 0.00261646873619519
$ws.Range("I25").Value = 1.965589387443117
$ws.Range("J25").Value = 0.09741669822220445
$ws.Range("K25").Value = 1.738742264864044
$ws.Range("L25").Value = 0.5892052198230999
$ws.Range("M25").Value = 0.3523585078118998
$ws.Range("N25").Value = 3.672429817011619

Write-Host "Applied 264 cell updates"